# Atualização automática de SAO_FRANCISCO_DE_ASSIS.xlsx
$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Keep the first sheet active/selected, matching the original workbook state
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
